$d = $word.ActiveDocument

# --- Step 1: modify the first paragraph's text & append colored runs ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Text = "This is a Microsoft word document.  "

# Collapse to the end of paragraph 1's text (before the paragraph mark)
$end1 = $r1.End - 1
$ins = $d.Range($end1, $end1)
$ins.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$ins.Font.Color = 192

$end1b = $end1 + ("(This is a change " + [char]0x2013 + " Ve").Length
$ins2 = $d.Range($end1b, $end1b)
$ins2.InsertAfter("rsion for branch alternate")
$ins2.Font.Color = 192

$end1c = $end1b + ("rsion for branch alternate").Length
$ins3 = $d.Range($end1c, $end1c)
$ins3.InsertAfter(")")
$ins3.Font.Color = 192

Write-Output ("Para1 now: [" + $d.Paragraphs.Item(1).Range.Text + "]")

# --- Step 2: insert a new empty paragraph after paragraph 2 ---
$p2 = $d.Paragraphs.Item(2)
$newPara = $p2.Range.InsertParagraphAfter()

